$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 1428686.1
$ws.Cells.Item(6, 9).Value = 1666717.1
$ws.Cells.Item(6, 11).Value = 5000151.300000001
$ws.Cells.Item(6, 13).Value = -5000039.300000001
$ws.Cells.Item(32, 8).Value = 2107.3333
$ws.Cells.Item(32, 9).Value = 1950
$ws.Cells.Item(32, 11).Value = 1950
$ws.Cells.Item(32, 13).Value = -1624
$ws.Cells.Item(40, 8).Value = 6527.75
$ws.Cells.Item(40, 9).Value = 3055.5
$ws.Cells.Item(40, 11).Value = 3055.5
$ws.Cells.Item(40, 13).Value = -2880.5
$ws.Cells.Item(51, 8).Value = 2871.4285
$ws.Cells.Item(51, 9).Value = 2557.6924
$ws.Cells.Item(51, 11).Value = 2557.6924
$ws.Cells.Item(51, 13).Value = -2073.6924
$ws.Cells.Item(92, 8).Value = 2230.158
$ws.Cells.Item(92, 9).Value = 2110.6924
$ws.Cells.Item(92, 10).Value = 2489
$ws.Cells.Item(92, 11).Value = 2110.6924
$ws.Cells.Item(92, 12).Value = 2489
$ws.Cells.Item(92, 13).Value = -862.6923999999999
$ws.Cells.Item(92, 14).Value = -4985
$ws.Cells.Item(137, 8).Value = 3932.55
$ws.Cells.Item(137, 10).Value = 2954.75
$ws.Cells.Item(137, 12).Value = 8864.25
$ws.Cells.Item(137, 14).Value = -13964.25
$ws.Cells.Item(138, 8).Value = 5578.47
$ws.Cells.Item(138, 10).Value = 5823.864
$ws.Cells.Item(138, 12).Value = 17471.592
$ws.Cells.Item(138, 14).Value = -27751.592

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 11842.75
$ws.Cells.Item(32, 9).Value = 10972.217
$ws.Cells.Item(32, 10).Value = 39990
$ws.Cells.Item(32, 11).Value = 10972.217
$ws.Cells.Item(32, 12).Value = 39990
$ws.Cells.Item(32, 13).Value = -10685.217
$ws.Cells.Item(32, 14).Value = -40564
$ws.Cells.Item(74, 8).Value = 44861.816
$ws.Cells.Item(74, 9).Value = 48941.715
$ws.Cells.Item(74, 11).Value = 48941.715
$ws.Cells.Item(74, 13).Value = -48067.715
$ws.Cells.Item(77, 8).Value = 44861.816
$ws.Cells.Item(77, 9).Value = 48941.715
$ws.Cells.Item(77, 11).Value = 244708.575
$ws.Cells.Item(77, 13).Value = -240340.575
$ws.Cells.Item(122, 8).Value = 3000
$ws.Cells.Item(122, 9).Value = 3000
$ws.Cells.Item(122, 11).Value = 9000
$ws.Cells.Item(122, 13).Value = -6550

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 2230.2
$ws.Cells.Item(86, 9).Value = 2251.75
$ws.Cells.Item(86, 11).Value = 2251.75
$ws.Cells.Item(86, 13).Value = -1128.75
$ws.Cells.Item(89, 8).Value = 2230.2
$ws.Cells.Item(89, 9).Value = 2251.75
$ws.Cells.Item(89, 11).Value = 11258.75
$ws.Cells.Item(89, 13).Value = -5642.75
$ws.Cells.Item(105, 8).Value = 1268.3636
$ws.Cells.Item(105, 9).Value = 1268.3636
$ws.Cells.Item(105, 11).Value = 1268.3636
$ws.Cells.Item(105, 13).Value = 478.6364000000001
$ws.Cells.Item(134, 8).Value = 3329.3333
$ws.Cells.Item(134, 9).Value = 2218.3635
$ws.Cells.Item(134, 10).Value = 5075.143
$ws.Cells.Item(134, 11).Value = 6655.0905
$ws.Cells.Item(134, 12).Value = 15225.429
$ws.Cells.Item(134, 13).Value = -4120.0905
$ws.Cells.Item(134, 14).Value = -20295.429

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 281.16666
$ws.Cells.Item(7, 9).Value = 30.555555
$ws.Cells.Item(7, 11).Value = 30.555555
$ws.Cells.Item(7, 13).Value = 82.444445
$ws.Cells.Item(31, 8).Value = 2843.7036
$ws.Cells.Item(31, 9).Value = 2574.3333
$ws.Cells.Item(31, 10).Value = 4998.6665
$ws.Cells.Item(31, 11).Value = 2574.3333
$ws.Cells.Item(31, 12).Value = 4998.6665
$ws.Cells.Item(31, 13).Value = -2279.3333
$ws.Cells.Item(31, 14).Value = -5588.6665
$ws.Cells.Item(34, 8).Value = 2843.7036
$ws.Cells.Item(34, 9).Value = 2574.3333
$ws.Cells.Item(34, 10).Value = 4998.6665
$ws.Cells.Item(34, 11).Value = 2574.3333
$ws.Cells.Item(34, 12).Value = 4998.6665
$ws.Cells.Item(34, 13).Value = -2372.3333
$ws.Cells.Item(34, 14).Value = -5402.6665
$ws.Cells.Item(58, 8).Value = 2739.88
$ws.Cells.Item(58, 9).Value = 2576.889
$ws.Cells.Item(58, 10).Value = 3159
$ws.Cells.Item(58, 11).Value = 2576.889
$ws.Cells.Item(58, 12).Value = 3159
$ws.Cells.Item(58, 13).Value = -2373.889
$ws.Cells.Item(58, 14).Value = -3565
$ws.Cells.Item(69, 8).Value = 10280
$ws.Cells.Item(69, 9).Value = 2336
$ws.Cells.Item(69, 11).Value = 2336
$ws.Cells.Item(69, 13).Value = -1587
$ws.Cells.Item(72, 8).Value = 10280
$ws.Cells.Item(72, 9).Value = 2336
$ws.Cells.Item(72, 11).Value = 7008
$ws.Cells.Item(72, 13).Value = -3264
$ws.Cells.Item(86, 8).Value = 3731.5
$ws.Cells.Item(86, 9).Value = 3643.4
$ws.Cells.Item(86, 11).Value = 3643.4
$ws.Cells.Item(86, 13).Value = -2520.4
$ws.Cells.Item(89, 8).Value = 3731.5
$ws.Cells.Item(89, 9).Value = 3643.4
$ws.Cells.Item(89, 11).Value = 18217
$ws.Cells.Item(89, 13).Value = -12601
$ws.Cells.Item(133, 8).Value = 98062.5
$ws.Cells.Item(133, 10).Value = 104928.57
$ws.Cells.Item(133, 12).Value = 104928.57
$ws.Cells.Item(133, 14).Value = -109988.57
$ws.Cells.Item(136, 8).Value = 2739.88
$ws.Cells.Item(136, 9).Value = 2576.889
$ws.Cells.Item(136, 10).Value = 3159
$ws.Cells.Item(136, 11).Value = 7730.667
$ws.Cells.Item(136, 12).Value = 9477
$ws.Cells.Item(136, 13).Value = -5180.667
$ws.Cells.Item(136, 14).Value = -14577

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(32, 8).Value = 1778.6666
$ws.Cells.Item(32, 9).Value = 168
$ws.Cells.Item(32, 11).Value = 504
$ws.Cells.Item(32, 13).Value = -221
$ws.Cells.Item(39, 8).Value = 1000
$ws.Cells.Item(39, 10).Value = 1000
$ws.Cells.Item(39, 12).Value = 3000
$ws.Cells.Item(39, 14).Value = -3588
$ws.Cells.Item(55, 8).Value = 1998.5714
$ws.Cells.Item(74, 8).Value = 10000
$ws.Cells.Item(74, 10).Value = 10000
$ws.Cells.Item(74, 12).Value = 30000
$ws.Cells.Item(74, 14).Value = -32122
$ws.Cells.Item(77, 8).Value = 10000
$ws.Cells.Item(77, 10).Value = 10000
$ws.Cells.Item(77, 12).Value = 90000
$ws.Cells.Item(77, 14).Value = -100608
$ws.Cells.Item(120, 8).Value = 966.6667
$ws.Cells.Item(120, 9).Value = 200
$ws.Cells.Item(120, 10).Value = 2500
$ws.Cells.Item(120, 11).Value = 600
$ws.Cells.Item(120, 12).Value = 7500
$ws.Cells.Item(120, 13).Value = 4238
$ws.Cells.Item(120, 14).Value = -17176
$ws.Cells.Item(131, 8).Value = 4272
$ws.Cells.Item(131, 10).Value = 4521.7393
$ws.Cells.Item(131, 12).Value = 13565.2179
$ws.Cells.Item(131, 14).Value = -23645.2179

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(41, 8).Value = 1114.6666
$ws.Cells.Item(41, 9).Value = 1114.6666
$ws.Cells.Item(41, 11).Value = 1114.6666
$ws.Cells.Item(41, 13).Value = -759.6666
$ws.Cells.Item(132, 8).Value = 3811.7273
$ws.Cells.Item(132, 9).Value = 2336.4614
$ws.Cells.Item(132, 11).Value = 7009.3842
$ws.Cells.Item(132, 13).Value = -4479.3842

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 4914.0713
$ws.Cells.Item(22, 9).Value = 4159.4
$ws.Cells.Item(22, 11).Value = 4159.4
$ws.Cells.Item(22, 13).Value = -3864.4
$ws.Cells.Item(27, 8).Value = 4914.0713
$ws.Cells.Item(27, 9).Value = 4159.4
$ws.Cells.Item(27, 11).Value = 4159.4
$ws.Cells.Item(27, 13).Value = -4052.4
$ws.Cells.Item(42, 8).Value = 25014
$ws.Cells.Item(42, 10).Value = 25028
$ws.Cells.Item(42, 12).Value = 25028
$ws.Cells.Item(42, 14).Value = -26154
$ws.Cells.Item(46, 8).Value = 2927.3572
$ws.Cells.Item(46, 10).Value = 3083.3076
$ws.Cells.Item(46, 12).Value = 3083.3076
$ws.Cells.Item(46, 14).Value = -3459.3076
$ws.Cells.Item(47, 8).Value = 0
$ws.Cells.Item(47, 10).Value = 0
$ws.Cells.Item(47, 12).Value = 0
$ws.Cells.Item(47, 14).ClearContents()  # N47 removed (was -31980)
$ws.Cells.Item(49, 8).Value = 25014
$ws.Cells.Item(49, 10).Value = 25028
$ws.Cells.Item(49, 12).Value = 25028
$ws.Cells.Item(49, 14).Value = -25322
$ws.Cells.Item(52, 8).Value = 0
$ws.Cells.Item(52, 10).Value = 0
$ws.Cells.Item(52, 12).Value = 0
$ws.Cells.Item(52, 14).ClearContents()  # N52 removed (was -31466)
$ws.Cells.Item(68, 8).Value = 3876.7693
$ws.Cells.Item(68, 9).Value = 3699.8333
$ws.Cells.Item(68, 11).Value = 3699.8333
$ws.Cells.Item(68, 13).Value = -2950.8333
$ws.Cells.Item(71, 8).Value = 3876.7693
$ws.Cells.Item(71, 9).Value = 3699.8333
$ws.Cells.Item(71, 11).Value = 18499.1665
$ws.Cells.Item(71, 13).Value = -14755.1665

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(46, 8).Value = 96250
$ws.Cells.Item(46, 10).Value = 96250
$ws.Cells.Item(46, 12).Value = 96250
$ws.Cells.Item(46, 14).Value = -96712
$ws.Cells.Item(62, 8).Value = 26999.75
$ws.Cells.Item(62, 9).Value = 24999.5
$ws.Cells.Item(62, 11).Value = 24999.5
$ws.Cells.Item(62, 13).Value = -24375.5
$ws.Cells.Item(65, 8).Value = 26999.75
$ws.Cells.Item(65, 9).Value = 24999.5
$ws.Cells.Item(65, 11).Value = 124997.5
$ws.Cells.Item(65, 13).Value = -121877.5
$ws.Cells.Item(122, 8).Value = 3399.25
$ws.Cells.Item(122, 9).Value = 3689.3
$ws.Cells.Item(122, 11).Value = 11067.9
$ws.Cells.Item(122, 13).Value = -8617.900000000001
$ws.Cells.Item(134, 8).Value = 96250
$ws.Cells.Item(134, 10).Value = 96250
$ws.Cells.Item(134, 12).Value = 288750
$ws.Cells.Item(134, 14).Value = -293820
